$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of cell updates from the refreshed "cryptos" price/volume feed.
# Each entry updates one cell to its new textual value. Columns D (Price) and
# E (Volume(1h)) are stored as plain text in the workbook (not numbers), so
# whenever the new text looks like a plain number (e.g. "4.33"), we force a
# text NumberFormat before assigning it so Excel keeps it as text instead of
# converting it to a numeric value, then restore the default "Normal" style
# once the text value has been written so no stray formatting is left behind.
$updates = @(
    @{ Cell = "D2"; Value = "62.856.76" },
    @{ Cell = "E2"; Value = "  -0.31%  " },
    @{ Cell = "D3"; Value = "2.544.67" },
    @{ Cell = "E3"; Value = "  +3.66%  " },
    @{ Cell = "E4"; Value = "  +0.03%  " },
    @{ Cell = "D5"; Value = "567.52" },
    @{ Cell = "E5"; Value = "  +0.78%  " },
    @{ Cell = "D6"; Value = "145.40" },
    @{ Cell = "E6"; Value = "  +2.30%  " },
    @{ Cell = "E7"; Value = "  +0.01%  " },
    @{ Cell = "D8"; Value = "0.581" },
    @{ Cell = "E8"; Value = "  +0.01%  " },
    @{ Cell = "D9"; Value = "2.542.53" },
    @{ Cell = "E9"; Value = "  +3.67%  " },
    @{ Cell = "D10"; Value = "0.104" },
    @{ Cell = "E10"; Value = "  -0.07%  " },
    @{ Cell = "D11"; Value = "5.49" },
    @{ Cell = "E11"; Value = "  -2.82%  " },
    @{ Cell = "E12"; Value = "  -0.05%  " },
    @{ Cell = "D13"; Value = "0.351" },
    @{ Cell = "E13"; Value = "  +0.06%  " },
    @{ Cell = "D14"; Value = "27.27" },
    @{ Cell = "E14"; Value = "  +1.91%  " },
    @{ Cell = "D15"; Value = "3.000.90" },
    @{ Cell = "E15"; Value = "  +3.72%  " },
    @{ Cell = "D16"; Value = "62.810.96" },
    @{ Cell = "E16"; Value = "  -0.15%  " },
    @{ Cell = "D17"; Value = "0.0000142" },
    @{ Cell = "E17"; Value = "  +1.45%  " },
    @{ Cell = "D18"; Value = "2.564.33" },
    @{ Cell = "E18"; Value = "  +4.56%  " },
    @{ Cell = "D19"; Value = "11.31" },
    @{ Cell = "E19"; Value = "  +0.95%  " },
    @{ Cell = "B20"; Value = "Polkadot" },
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot" },
    @{ Cell = "D20"; Value = "4.33" },
    @{ Cell = "E20"; Value = "  +2.05%  " },
    @{ Cell = "B21"; Value = "BitcoinCash" },
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch" },
    @{ Cell = "D21"; Value = "332.39" },
    @{ Cell = "E21"; Value = "  -1.74%  " },
    @{ Cell = "D22"; Value = "6.80" },
    @{ Cell = "E22"; Value = "  +1.29%  " },
    @{ Cell = "E23"; Value = "  +0.08%  " },
    @{ Cell = "D24"; Value = "65.16" },
    @{ Cell = "E24"; Value = "  +0.06%  " },
    @{ Cell = "D25"; Value = "0.169" },
    @{ Cell = "E25"; Value = "  -0.98%  " },
    @{ Cell = "D26"; Value = "1.59" },
    @{ Cell = "E26"; Value = "  +7.29%  " },
    @{ Cell = "E27"; Value = "  +0.18%  " },
    @{ Cell = "D28"; Value = "8.37" },
    @{ Cell = "E28"; Value = "  +4.76%  " },
    @{ Cell = "D29"; Value = "1.46" },
    @{ Cell = "E29"; Value = "  +4.08%  " },
    @{ Cell = "D30"; Value = "7.29" },
    @{ Cell = "E30"; Value = "  +7.56%  " },
    @{ Cell = "D31"; Value = "0.0₃0810" },
    @{ Cell = "E31"; Value = "  +3.05%  " },
    @{ Cell = "D32"; Value = "1.84" },
    @{ Cell = "E32"; Value = "  +0.68%  " },
    @{ Cell = "D33"; Value = "174.61" },
    @{ Cell = "E33"; Value = "  -1.41%  " },
    @{ Cell = "D34"; Value = "1.55" },
    @{ Cell = "E34"; Value = "  +2.91%  " },
    @{ Cell = "D35"; Value = "405.70" },
    @{ Cell = "E35"; Value = "  +5.98%  " },
    @{ Cell = "D36"; Value = "0.398" },
    @{ Cell = "E36"; Value = "  +0.49%  " },
    @{ Cell = "D37"; Value = "18.88" },
    @{ Cell = "E37"; Value = "  +0.85%  " },
    @{ Cell = "E38"; Value = "  +0.00%  " },
    @{ Cell = "D39"; Value = "4.33" },
    @{ Cell = "E39"; Value = "  +0.45%  " },
    @{ Cell = "D40"; Value = "1.73" },
    @{ Cell = "E40"; Value = "  +1.50%  " },
    @{ Cell = "E41"; Value = "  +0.06%  " },
    @{ Cell = "D42"; Value = "39.61" },
    @{ Cell = "E42"; Value = "  -1.05%  " },
    @{ Cell = "D43"; Value = "151.53" },
    @{ Cell = "E43"; Value = "  +1.85%  " },
    @{ Cell = "D44"; Value = "3.75" },
    @{ Cell = "E44"; Value = "  +1.63%  " },
    @{ Cell = "D45"; Value = "20.64" },
    @{ Cell = "E45"; Value = "  +1.42%  " },
    @{ Cell = "D46"; Value = "0.603" },
    @{ Cell = "E46"; Value = "  +1.36%  " },
    @{ Cell = "B47"; Value = "Hedera" },
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar" },
    @{ Cell = "D47"; Value = "0.0528" },
    @{ Cell = "E47"; Value = "  +2.96%  " },
    @{ Cell = "B48"; Value = "Stellar" },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm" },
    @{ Cell = "D48"; Value = "0.0961" },
    @{ Cell = "E48"; Value = "  +0.21%  " },
    @{ Cell = "B49"; Value = "VeChain" },
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet" },
    @{ Cell = "D49"; Value = "0.0237" },
    @{ Cell = "E49"; Value = "  +4.02%  " },
    @{ Cell = "D50"; Value = "18.24" },
    @{ Cell = "E50"; Value = "  +2.60%  " },
    @{ Cell = "D51"; Value = "1.74" },
    @{ Cell = "E51"; Value = "  -0.82%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Value -match '^[+-]?\d+(\.\d+)?$') {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
